$wb = $excel.ActiveWorkbook

# "Generate Report for Handback": refresh the recorded handoff/handback
# timestamps for the 569d02c7-... entry (row 12) on the zh-cn and de-de
# status sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D12").Value = "2016-03-03 07:41:05"
$wsZhCn.Range("G12").Value = "2016-03-03 07:41:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D12").Value = "2016-03-03 07:41:16"
$wsDeDe.Range("G12").Value = "2016-03-03 07:42:10"
